$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '43.210.03'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +1.87%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.382.35'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +4.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '303.11'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.74%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '96.98'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  -0.09%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.501'
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '34.21'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.18%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0789'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.17%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.121'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +2.42%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '18.41'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.18%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.79'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.19%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '2.755.10'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +4.02%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.374.94'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +3.34%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.809'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +4.10%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '43.225.01'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.02%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.19'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.33'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +6.41%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0₃0889'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.31%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '68.68'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.84%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '235.11'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('E25').Value = '  -0.23%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.43'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.16%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '24.83'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +2.69%  '
$ws.Range('E28').Value = '  +3.03%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.13'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.24%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '31.53'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  -0.05%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '5.10'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +2.72%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0735'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +5.99%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '17.09'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.24%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.85'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +7.21%  '
$ws.Range('E36').Value = '  +2.70%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.32'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.30'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.01%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.79'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +4.82%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '22.23'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +11.76%  '
$ws.Range('E41').Value = '  +0.34%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '105.72'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -35.96%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.955.00'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.26%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0280'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  +1.71%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.88%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '9.26'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -10.29%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '52.79'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('E49').Value = '  +3.46%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '71.91'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.08%  '
$ws.Range('E51').Value = '  +1.45%  '
